# MAuS Checklist.xlsx - "display done or not done in TodoListActivity (63/100)"
#
# Row 28 ("das Erledigsein/Nicht-Erledigtsein") is marked as implemented/done:
# the old "Sollte schnell machbar sein" note in E28 is cleared and D28 (the
# "done" column) is set to 1, same as the other completed requirement rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E28").Clear()
$ws.Range("D28").Value2 = 1

# Move the on-screen selection the same way the author's Excel session ended
# up (cosmetic, but reflected in the saved sheetView/selection).
$null = $ws.Range("E28").Select()

Write-Host "Row 28 marked done; D55 total will recalc automatically."
